$wb = $excel.ActiveWorkbook

# This script applies updated market-price derived values (columns H-N)
# across several worksheets, as produced by the scheduled pricing runner.

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1248.6
$ws.Range("I28").Value = 1248.5264
$ws.Range("K28").Value = 1248.5264
$ws.Range("M28").Value = -763.5264
$ws.Range("H33").Value = 2062.8696
$ws.Range("I33").Value = 2253
$ws.Range("J33").Value = 1159.75
$ws.Range("K33").Value = 2253
$ws.Range("L33").Value = 1159.75
$ws.Range("M33").Value = -2024
$ws.Range("N33").Value = -1617.75
$ws.Range("H40").Value = 3407.9092
$ws.Range("I40").Value = 2127.1428
$ws.Range("J40").Value = 5649.25
$ws.Range("K40").Value = 2127.1428
$ws.Range("L40").Value = 5649.25
$ws.Range("M40").Value = -1952.1428
$ws.Range("N40").Value = -5999.25
$ws.Range("H97").Value = 1806.2307
$ws.Range("J97").Value = 1873.4166
$ws.Range("L97").Value = 5620.2498
$ws.Range("N97").Value = -6612.2498
$ws.Range("H101").Value = 489.57144
$ws.Range("I101").Value = 481.75
$ws.Range("K101").Value = 1445.25
$ws.Range("M101").Value = 176.75
$ws.Range("H107").Value = 823.8261
$ws.Range("I107").Value = 515.3077
$ws.Range("J107").Value = 1224.9
$ws.Range("K107").Value = 515.3077
$ws.Range("L107").Value = 1224.9
$ws.Range("M107").Value = 1404.6923
$ws.Range("N107").Value = -5064.9
$ws.Range("H132").Value = 2379.7793
$ws.Range("I132").Value = 1860.6666
$ws.Range("K132").Value = 5581.9998
$ws.Range("M132").Value = -3051.9998
$ws.Range("H138").Value = 4055.5789
$ws.Range("J138").Value = 4098.7646
$ws.Range("L138").Value = 12296.2938
$ws.Range("N138").Value = -22576.2938

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8212.6
$ws.Range("I32").Value = 5505.321
$ws.Range("J32").Value = 19754.158
$ws.Range("K32").Value = 5505.321
$ws.Range("L32").Value = 19754.158
$ws.Range("M32").Value = -5218.321
$ws.Range("N32").Value = -20328.158
$ws.Range("H63").Value = 6070
$ws.Range("I63").Value = 6070
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 6070
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -5384
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 6070
$ws.Range("I66").Value = 6070
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 30350
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -26918
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 2840.2334
$ws.Range("I74").Value = 2495.4583
$ws.Range("J74").Value = 4219.3335
$ws.Range("K74").Value = 2495.4583
$ws.Range("L74").Value = 4219.3335
$ws.Range("M74").Value = -1621.4583
$ws.Range("N74").Value = -5967.3335
$ws.Range("H77").Value = 2840.2334
$ws.Range("I77").Value = 2495.4583
$ws.Range("J77").Value = 4219.3335
$ws.Range("K77").Value = 12477.2915
$ws.Range("L77").Value = 21096.6675
$ws.Range("M77").Value = -8109.291499999999
$ws.Range("N77").Value = -29832.6675
$ws.Range("H108").Value = 49999.668
$ws.Range("J108").Value = 49999.668
$ws.Range("L108").Value = 49999.668
$ws.Range("N108").Value = -57679.668
$ws.Range("H122").Value = 5616.727
$ws.Range("J122").Value = 6734.143
$ws.Range("L122").Value = 20202.429
$ws.Range("N122").Value = -25102.429
$ws.Range("H132").Value = 7910.655
$ws.Range("I132").Value = 8226.291999999999
$ws.Range("J132").Value = 5746.2856
$ws.Range("K132").Value = 24678.876
$ws.Range("L132").Value = 17238.8568
$ws.Range("M132").Value = -22148.876
$ws.Range("N132").Value = -22298.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 256.8
$ws.Range("I22").Value = 256.8
$ws.Range("K22").Value = 256.8
$ws.Range("M22").Value = -83.80000000000001
$ws.Range("H45").Value = 34999
$ws.Range("I45").Value = 34999
$ws.Range("K45").Value = 34999
$ws.Range("M45").Value = -34191
$ws.Range("H99").Value = 3462.7778
$ws.Range("I99").Value = 2786
$ws.Range("K99").Value = 2786
$ws.Range("M99").Value = -1288

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 472
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 577.5
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 577.5
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = -801.5
$ws.Range("H14").Value = 1599.5
$ws.Range("I14").Value = 199
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 199
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -29
$ws.Range("N14").Value = -3340
$ws.Range("H15").Value = 25708.5
$ws.Range("I15").Value = 9999
$ws.Range("J15").Value = 30945
$ws.Range("K15").Value = 9999
$ws.Range("L15").Value = 30945
$ws.Range("M15").Value = -9829
$ws.Range("N15").Value = -31285
$ws.Range("H16").Value = 2045
$ws.Range("I16").Value = 817.5714
$ws.Range("J16").Value = 2904.2
$ws.Range("K16").Value = 817.5714
$ws.Range("L16").Value = 2904.2
$ws.Range("M16").Value = -530.5714
$ws.Range("N16").Value = -3478.2
$ws.Range("H23").Value = 1986639.8
$ws.Range("J23").Value = 1986639.8
$ws.Range("L23").Value = 1986639.8
$ws.Range("N23").Value = -1987119.8
$ws.Range("H27").Value = 1986639.8
$ws.Range("J27").Value = 1986639.8
$ws.Range("L27").Value = 1986639.8
$ws.Range("N27").Value = -1987023.8
$ws.Range("H58").Value = 3361.9666
$ws.Range("I58").Value = 1495.2941
$ws.Range("J58").Value = 5803
$ws.Range("K58").Value = 1495.2941
$ws.Range("L58").Value = 5803
$ws.Range("M58").Value = -1292.2941
$ws.Range("N58").Value = -6209
$ws.Range("H109").Value = 34999.5
$ws.Range("J109").Value = 34999.5
$ws.Range("L109").Value = 34999.5
$ws.Range("N109").Value = -37079.5
$ws.Range("H113").Value = 2045
$ws.Range("I113").Value = 817.5714
$ws.Range("J113").Value = 2904.2
$ws.Range("K113").Value = 817.5714
$ws.Range("L113").Value = 2904.2
$ws.Range("M113").Value = 1352.4286
$ws.Range("N113").Value = -7244.2
$ws.Range("H132").Value = 5059.863
$ws.Range("I132").Value = 3021.3547
$ws.Range("J132").Value = 16549.637
$ws.Range("K132").Value = 9064.0641
$ws.Range("L132").Value = 49648.91099999999
$ws.Range("M132").Value = -6534.0641
$ws.Range("N132").Value = -54708.91099999999
$ws.Range("H136").Value = 3361.9666
$ws.Range("I136").Value = 1495.2941
$ws.Range("J136").Value = 5803
$ws.Range("K136").Value = 4485.8823
$ws.Range("L136").Value = 17409
$ws.Range("M136").Value = -1935.8823
$ws.Range("N136").Value = -22509
$ws.Range("H141").Value = 238000.2
$ws.Range("J141").Value = 238000.2
$ws.Range("L141").Value = 238000.2
$ws.Range("N141").Value = -248360.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 8848.4
$ws.Range("I57").Value = 8848.4
$ws.Range("K57").Value = 26545.2
$ws.Range("M57").Value = -25986.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12413.235
$ws.Range("I80").Value = 2868.2
$ws.Range("K80").Value = 2868.2
$ws.Range("M80").Value = -1870.2
$ws.Range("H83").Value = 12413.235
$ws.Range("I83").Value = 2868.2
$ws.Range("K83").Value = 14341
$ws.Range("M83").Value = -9349
$ws.Range("H102").Value = 41406.027
$ws.Range("I102").Value = 55012.883
$ws.Range("K102").Value = 55012.883
$ws.Range("M102").Value = -53390.883
$ws.Range("H126").Value = 66597.48
$ws.Range("I126").Value = 95253.69500000001
$ws.Range("K126").Value = 285761.085
$ws.Range("M126").Value = -283291.085
$ws.Range("H132").Value = 5906.091
$ws.Range("I132").Value = 4937.353
$ws.Range("J132").Value = 9199.799999999999
$ws.Range("K132").Value = 14812.059
$ws.Range("L132").Value = 27599.4
$ws.Range("M132").Value = -12282.059
$ws.Range("N132").Value = -32659.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 29872.857
$ws.Range("J104").Value = 29872.857
$ws.Range("L104").Value = 29872.857
$ws.Range("N104").Value = -36860.857
